# Subjects_Concepts_Links.xlsx — "DBMS & CN Links updated"
#
# 1. Add two new sheets "OS" and "DBMS" after "CN".
# 2. Populate OS with a couple of reference rows.
# 3. Populate DBMS with just the header row.
# 4. Add a new row (LAN / Token Ring) to the CN sheet.
# 5. Leave CN as the active/selected sheet (matches original activeTab).

$wb = $excel.ActiveWorkbook

$cn = $wb.Worksheets.Item("CN")

# --- add the OS sheet right after CN ---------------------------------
$os = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $cn)
$os.Name = "OS"

# --- add the DBMS sheet right after OS --------------------------------
$dbms = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $os)
$dbms.Name = "DBMS"

# -----------------------------------------------------------------
# OS sheet content
# Header first
$os.Range("A1").Value = "Chapter"
$os.Range("B1").Value = "Topic"
$os.Range("C1").Value = "Link"
$os.Range("A1:C1").Font.Bold = $true

# Row 6 (Process Management / Scheduling Algorithms) -- written before
# row 3 so new shared strings land in the same order as the source file.
$os.Range("A6").Value = "Process Management"
$os.Range("B6").Value = "Schedling Algorithms"
$os.Range("C6").Value = "https://www.youtube.com/playlist?list=PLEbnTDJUr_If_BnzJkkN_J0Tl3iXTL8vq"

# Row 3 (Reading Material / Article from JavaTpoint)
$os.Range("A3").Value = "Reading Material"
$os.Range("B3").Value = "Article from JavaTpoint"
$os.Range("C3").Value = "https://www.javatpoint.com/os-tutorial"

$os.Columns.Item(1).ColumnWidth = 17.833333333333332
$os.Columns.Item(2).ColumnWidth = 19.333333333333332
$os.Columns.Item(3).ColumnWidth = 65.16666666666667

$os.Range("C3").Select() | Out-Null

# -----------------------------------------------------------------
# DBMS sheet content -- header row only
$dbms.Range("A1").Value = "Chapter"
$dbms.Range("B1").Value = "Topic"
$dbms.Range("C1").Value = "Link"
$dbms.Range("A1:C1").Font.Bold = $true

$dbms.Rows.Item(1).Select() | Out-Null

# -----------------------------------------------------------------
# CN sheet gets a new row of links (LAN / Token Ring)
$cn.Range("A25").Value = "LAN"
$cn.Range("C25").Value = "http://www.cs.montana.edu/~halla/csci466/lectures/lec10-2.7-token.html"
$cn.Range("B25").Value = "Token Ring (**Priority Scheme)"

$cn.Columns.Item(2).ColumnWidth = 25.666666666666668

# Restore CN as the active sheet / selection so the workbook opens where
# it did before (CN was already the active tab).
$cn.Activate()
$cn.Range("B25").Select() | Out-Null

Write-Host "OS/DBMS sheets added, CN updated"
